$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column K / L notes about data types (rows 11-23)
# Order matters for shared-string allocation, so set them in the same
# order the original author typed them in.
$ws.Range("K11").Value = "Categorical data"
$ws.Range("K12").Value = "Binary data"
$ws.Range("L13").Value = "Yes/No"
$ws.Range("L14").Value = "Dead/Alive"
$ws.Range("K15").Value = "Nominal (”label”, several groups)"
$ws.Range("L16").Value = "Eye colour: Blue/ Brown / Grey / Green"
$ws.Range("L17").Value = "Where do you live: Denmark, Germany, Sweden."
$ws.Range("K18").Value = "Ordinal"
$ws.Range("L19").Value = "How do you feel today?: Very unhappy, unhappy, OK, happy, very"
$ws.Range("L20").Value = "happy."
$ws.Range("L21").Value = "Do you try to eat healthily?: Never, Sometimes, Always"
$ws.Range("K22").Value = "Interval (does have a numerical distance between values)"
$ws.Range("L23").Value = "BMI categories (<25, 25-"

# Row 46 additions (extra example row for PCA prep)
$ws.Range("E46").Value = 1
$ws.Range("F46").Value = "numeric"
$ws.Range("G46").Value = 3
$ws.Range("H46").Value = "all 3 numeric"

# New header label for column I (row 1), added after the row 46 notes
$ws.Range("I1").Value = "other notes"
$ws.Range("I46").Value = "linear regrssion to simple"

# Update selection / view to match author's final state
$ws.Range("K11:L23").Select()
